$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 700.1667
$ws.Range("I6").Value = 334.5
$ws.Range("K6").Value = 1003.5
$ws.Range("M6").Value = -891.5

$ws.Range("H19").Value = 4683.76
$ws.Range("I19").Value = 4457.8667
$ws.Range("J19").Value = 5022.6
$ws.Range("K19").Value = 4457.8667
$ws.Range("L19").Value = 5022.6
$ws.Range("M19").Value = -4282.8667
$ws.Range("N19").Value = -5372.6

$ws.Range("H64").Value = 54500.5
$ws.Range("J64").Value = 6000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6496

$ws.Range("H67").Value = 54500.5
$ws.Range("J67").Value = 6000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7716

$ws.Range("H80").Value = 15267.6
$ws.Range("J80").Value = 13731.846
$ws.Range("L80").Value = 41195.538
$ws.Range("N80").Value = -43191.538

$ws.Range("H83").Value = 15267.6
$ws.Range("J83").Value = 13731.846
$ws.Range("L83").Value = 123586.614
$ws.Range("N83").Value = -133570.614

$ws.Range("H86").Value = 3257.8
$ws.Range("I86").Value = 3674.25
$ws.Range("J86").Value = 2980.1667
$ws.Range("K86").Value = 3674.25
$ws.Range("L86").Value = 2980.1667
$ws.Range("M86").Value = -2551.25
$ws.Range("N86").Value = -5226.1667

$ws.Range("H89").Value = 3257.8
$ws.Range("I89").Value = 3674.25
$ws.Range("J89").Value = 2980.1667
$ws.Range("K89").Value = 18371.25
$ws.Range("L89").Value = 14900.8335
$ws.Range("M89").Value = -12755.25
$ws.Range("N89").Value = -26132.8335

$ws.Range("H99").Value = 566.8889
$ws.Range("I99").Value = 498.25
$ws.Range("J99").Value = 621.8
$ws.Range("K99").Value = 1494.75
$ws.Range("L99").Value = 1865.4
$ws.Range("M99").Value = 3.25
$ws.Range("N99").Value = -4861.4

$ws.Range("H106").Value = 1594.875
$ws.Range("I106").Value = 1500
$ws.Range("J106").Value = 1651.8
$ws.Range("K106").Value = 1500
$ws.Range("L106").Value = 1651.8
$ws.Range("M106").Value = -869
$ws.Range("N106").Value = -2913.8

$ws.Range("H129").Value = 394.9091
$ws.Range("I129").Value = 394.9091
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1184.7273
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3815.2727
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 14968.096
$ws.Range("I132").Value = 16604.043
$ws.Range("K132").Value = 49812.129
$ws.Range("M132").Value = -47282.129

$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17791.89
$ws.Range("I32").Value = 17448
$ws.Range("K32").Value = 17448
$ws.Range("M32").Value = -17161

$ws.Range("H63").Value = 4179.1875
$ws.Range("I63").Value = 2983.8572
$ws.Range("J63").Value = 5108.8887
$ws.Range("K63").Value = 2983.8572
$ws.Range("L63").Value = 5108.8887
$ws.Range("M63").Value = -2297.8572
$ws.Range("N63").Value = -6480.8887

$ws.Range("H66").Value = 4179.1875
$ws.Range("I66").Value = 2983.8572
$ws.Range("J66").Value = 5108.8887
$ws.Range("K66").Value = 14919.286
$ws.Range("L66").Value = 25544.4435
$ws.Range("M66").Value = -11487.286
$ws.Range("N66").Value = -32408.4435

$ws.Range("H110").Value = 3092.4
$ws.Range("I110").Value = 3078.6956
$ws.Range("J110").Value = 3250
$ws.Range("K110").Value = 3078.6956
$ws.Range("L110").Value = 3250
$ws.Range("M110").Value = -1033.6956
$ws.Range("N110").Value = -7340

$ws.Range("H122").Value = 2290
$ws.Range("I122").Value = 2311.9412
$ws.Range("K122").Value = 6935.823600000001
$ws.Range("M122").Value = -4485.823600000001

$ws.Range("H132").Value = 1356.3793
$ws.Range("I132").Value = 921.52
$ws.Range("J132").Value = 4074.25
$ws.Range("K132").Value = 2764.56
$ws.Range("L132").Value = 12222.75
$ws.Range("M132").Value = -234.5599999999999
$ws.Range("N132").Value = -17282.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 21765.885
$ws.Range("I107").Value = 24814.455
$ws.Range("J107").Value = 4998.75
$ws.Range("K107").Value = 24814.455
$ws.Range("L107").Value = 4998.75
$ws.Range("M107").Value = -22894.455
$ws.Range("N107").Value = -8838.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 478.26666
$ws.Range("I107").Value = 310.78262
$ws.Range("K107").Value = 310.78262
$ws.Range("M107").Value = 1609.21738

$ws.Range("H134").Value = 2130.2554
$ws.Range("I134").Value = 1909.814
$ws.Range("K134").Value = 5729.442
$ws.Range("M134").Value = -3194.442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 20299.666
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 20299.666
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 60898.99800000001
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -61066.99800000001

$ws.Range("H80").Value = 4967.72
$ws.Range("J80").Value = 4966.375
$ws.Range("L80").Value = 14899.125
$ws.Range("N80").Value = -16771.125

$ws.Range("H83").Value = 4967.72
$ws.Range("J83").Value = 4966.375
$ws.Range("L83").Value = 44697.375
$ws.Range("N83").Value = -54057.375

$ws.Range("H138").Value = 7334.1665
$ws.Range("I138").Value = 9743.091
$ws.Range("J138").Value = 3548.7144
$ws.Range("K138").Value = 29229.273
$ws.Range("L138").Value = 10646.1432
$ws.Range("M138").Value = -24089.273
$ws.Range("N138").Value = -20926.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 4286314.5
$ws.Range("I10").Value = 6000080
$ws.Range("K10").Value = 6000080
$ws.Range("M10").Value = -5999911

$ws.Range("H102").Value = 2454.3333
$ws.Range("I102").Value = 2150.8262
$ws.Range("J102").Value = 4199.5
$ws.Range("K102").Value = 2150.8262
$ws.Range("L102").Value = 4199.5
$ws.Range("M102").Value = -528.8262
$ws.Range("N102").Value = -7443.5

$ws.Range("H122").Value = 3863.1
$ws.Range("I122").Value = 3028.0435
$ws.Range("K122").Value = 9084.130500000001
$ws.Range("M122").Value = -6634.130500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10284

$ws.Range("H62").Value = 6777.5
$ws.Range("J62").Value = 8740.4
$ws.Range("L62").Value = 8740.4
$ws.Range("N62").Value = -9988.4

$ws.Range("H65").Value = 6777.5
$ws.Range("J65").Value = 8740.4
$ws.Range("L65").Value = 43702
$ws.Range("N65").Value = -49942

$ws.Range("H122").Value = 71225.266
$ws.Range("I122").Value = 83142.75
$ws.Range("K122").Value = 249428.25
$ws.Range("M122").Value = -246978.25
